# Update "想去人数" (want-to-go count) figures in the 展览, 演出 and 全部类型 sheets
# to reflect the newly generated gh-pages output (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (exhibitions) sheet -------------------------------------------------
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value  = 6981
$wsExhibit.Range("F4").Value  = 457
$wsExhibit.Range("F7").Value  = 129
$wsExhibit.Range("F15").Value = 1819
$wsExhibit.Range("F16").Value = 41
$wsExhibit.Range("F17").Value = 3591
$wsExhibit.Range("F19").Value = 244
$wsExhibit.Range("F21").Value = 16
$wsExhibit.Range("F23").Value = 2208
$wsExhibit.Range("F24").Value = 4
$wsExhibit.Range("F25").Value = 235
$wsExhibit.Range("F26").Value = 10
$wsExhibit.Range("F31").Value = 152
$wsExhibit.Range("F32").Value = 136
$wsExhibit.Range("F33").Value = 65

# --- 演出 (performances) sheet ------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F2").Value = 6

# --- 全部类型 (all types) sheet -----------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 6981
$wsAll.Range("F4").Value  = 457
$wsAll.Range("F6").Value  = 6
$wsAll.Range("F8").Value  = 129
$wsAll.Range("F16").Value = 1819
$wsAll.Range("F17").Value = 41
$wsAll.Range("F18").Value = 3591
$wsAll.Range("F20").Value = 244
$wsAll.Range("F22").Value = 16
$wsAll.Range("F24").Value = 2208
$wsAll.Range("F25").Value = 4
$wsAll.Range("F26").Value = 235
$wsAll.Range("F27").Value = 10
$wsAll.Range("F32").Value = 152
$wsAll.Range("F33").Value = 136
$wsAll.Range("F34").Value = 65

$wb.Save()
